$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$ws.Range("E9").Value = "INSTALA UNA RED LAN"
$ws.Range("E10").Value = "OPERA UNA RED LAN"
$ws.Range("E13").Value = "OPERA UNA RED LAN"
$ws.Range("E14").Value = "INSTALA UNA RED LAN"
